# Pio's ERA operativo Abril-Diciembre 2025
# Update computed scheduling metrics (seniority, theoretical/assigned load, shift
# counts, and per-month averages) on Sheet1 to reflect the refreshed "as of" date
# used for recalculation. These are externally computed cached values (no live
# formulas in the sheet), so each changed cell is written directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 11.3013698630137
$ws.Range("F2").Value = 6.221577520305049
$ws.Range("G2").Value = 4.321250888415055
$ws.Range("H2").Value = 22
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 1
$ws.Range("O2").Value = 5.527272727272726
$ws.Range("S2").Value = 0.3349862258953167
$ws.Range("C3").Value = 11.3013698630137
$ws.Range("F3").Value = 6.221577520305049
$ws.Range("G3").Value = 3.894811656005675
$ws.Range("H3").Value = 19
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 4.734426229508196
$ws.Range("S3").Value = 0.3322404371584699
$ws.Range("C4").Value = 10.7972602739726
$ws.Range("F4").Value = 6.277886928260588
$ws.Range("G4").Value = 4.321250888415055
$ws.Range("H4").Value = 22
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 5.266141732283463
$ws.Range("S4").Value = 0.3191601049868766
$ws.Range("C5").Value = 10.7972602739726
$ws.Range("F5").Value = 6.277886928260588
$ws.Range("G5").Value = 4.179104477611928
$ws.Range("H5").Value = 21
$ws.Range("M5").Value = 1
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 5.1072
$ws.Range("S5").Value = 0.3242666666666666
$ws.Range("C6").Value = 10.21643835616438
$ws.Range("F6").Value = 6.342765159165885
$ws.Range("G6").Value = 3.894811656005675
$ws.Range("H6").Value = 19
$ws.Range("M6").Value = 1
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 4.658064516129031
$ws.Range("S6").Value = 0.3268817204301075
$ws.Range("C7").Value = 10.21643835616438
$ws.Range("F7").Value = 6.342765159165885
$ws.Range("G7").Value = 4.434968017057556
$ws.Range("I7").Value = 19
$ws.Range("M7").Value = 1
$ws.Range("N7").Value = 1
$ws.Range("P7").Value = 4.548031496062992
$ws.Range("S7").Value = 0.3191601049868766
$ws.Range("C8").Value = 10.21643835616438
$ws.Range("F8").Value = 6.342765159165885
$ws.Range("G8").Value = 4.434968017057556
$ws.Range("I8").Value = 19
$ws.Range("M8").Value = 1
$ws.Range("N8").Value = 1
$ws.Range("P8").Value = 4.658064516129031
$ws.Range("S8").Value = 0.3268817204301075
$ws.Range("C9").Value = 7.794520547945205
$ws.Range("F9").Value = 6.613295140865327
$ws.Range("G9").Value = 8.955223880596993
$ws.Range("I9").Value = 22
$ws.Range("K9").Value = 12
$ws.Range("M9").Value = 1
$ws.Range("N9").Value = 1
$ws.Range("P9").Value = 5.393548387096773
$ws.Range("R9").Value = 2.941935483870967
$ws.Range("S9").Value = 0.4903225806451612
$ws.Range("C10").Value = 6.961643835616439
$ws.Range("F10").Value = 6.706328075748393
$ws.Range("G10").Value = 8.955223880596993
$ws.Range("I10").Value = 22
$ws.Range("K10").Value = 12
$ws.Range("M10").Value = 1
$ws.Range("N10").Value = 1
$ws.Range("P10").Value = 4.811510791366906
$ws.Range("R10").Value = 2.624460431654676
$ws.Range("S10").Value = 0.437410071942446
$ws.Range("C11").Value = 6.761643835616439
$ws.Range("F11").Value = 6.728668221295972
$ws.Range("G11").Value = 8.514570007107297
$ws.Range("J11").Value = 3
$ws.Range("Q11").Value = 0.7475409836065573
$ws.Range("R11").Value = 3.239344262295082
$ws.Range("S11").Value = 1.162841530054645
$ws.Range("C12").Value = 4.961643835616439
$ws.Range("F12").Value = 6.929729531224177
$ws.Range("G12").Value = 8.514570007107297
$ws.Range("J12").Value = 3
$ws.Range("K12").Value = 13
$ws.Range("Q12").Value = 0.7295999999999999
$ws.Range("R12").Value = 3.1616
$ws.Range("S12").Value = 1.134933333333333
$ws.Range("C13").Value = 3.126027397260274
$ws.Range("F13").Value = 7.134769223236198
$ws.Range("G13").Value = 8.727789623311988
$ws.Range("K13").Value = 13
$ws.Range("Q13").Value = 0.9806451612903224
$ws.Range("R13").Value = 3.187096774193548
$ws.Range("S13").Value = 1.144086021505376
$ws.Range("C14").Value = 2.378082191780822
$ws.Range("F14").Value = 7.218315246996319
$ws.Range("G14").Value = 8.727789623311988
$ws.Range("K14").Value = 13
$ws.Range("L14").Value = 4.666666666666666
$ws.Range("N14").Value = 2
$ws.Range("Q14").Value = 0.9727999999999999
$ws.Range("R14").Value = 3.1616
$ws.Range("S14").Value = 1.134933333333333
$ws.Range("C15").Value = 2.378082191780822
$ws.Range("F15").Value = 7.218315246996319
$ws.Range("G15").Value = 9.026297085998554
$ws.Range("K15").Value = 13
$ws.Range("L15").Value = 4.999999999999999
$ws.Range("M15").Value = 2
$ws.Range("Q15").Value = 0.874820143884892
$ws.Range("R15").Value = 2.843165467625899
$ws.Range("S15").Value = 1.093525179856115
$ws.Range("C16").Value = 0.5424657534246575
$ws.Range("F16").Value = 7.42335493900834
$ws.Range("G16").Value = 9.097370291400118
$ws.Range("K16").Value = 14
$ws.Range("L16").Value = 4.999999999999999
$ws.Range("N16").Value = 2
$ws.Range("Q16").Value = 0.6608695652173913
$ws.Range("R16").Value = 3.084057971014492
$ws.Range("S16").Value = 1.101449275362319
